# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# Update the conversion summary text on Hoja1!A1
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.6 = 9695.95 pesos`n✅ 9695.95 pesos = 2.59 = 924.25 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# Update the rate figures on the "tasas" sheet
$wsTasas.Range("N10").Value = 384.8
$wsTasas.Range("O10").Value = 3731
$wsTasas.Range("N12").Value = 3746
$wsTasas.Range("O12").Value = 357.081
